$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1194
$ws.Range("I32").Value = 785.3333
$ws.Range("J32").Value = 1500.5
$ws.Range("K32").Value = 785.3333
$ws.Range("L32").Value = 1500.5
$ws.Range("M32").Value = -459.3333
$ws.Range("N32").Value = -2152.5

$ws.Range("H62").Value = 3100.7778
$ws.Range("I62").Value = 2319
$ws.Range("J62").Value = 4078
$ws.Range("K62").Value = 2319
$ws.Range("L62").Value = 4078
$ws.Range("M62").Value = -1695
$ws.Range("N62").Value = -5326

$ws.Range("H65").Value = 3100.7778
$ws.Range("I65").Value = 2319
$ws.Range("J65").Value = 4078
$ws.Range("K65").Value = 11595
$ws.Range("L65").Value = 20390
$ws.Range("M65").Value = -8475
$ws.Range("N65").Value = -26630

$ws.Range("H93").Value = 50601
$ws.Range("J93").Value = 50601
$ws.Range("L93").Value = 50601
$ws.Range("N93").Value = -55593

$ws.Range("H129").Value = 921.6799999999999
$ws.Range("J129").Value = 1092.4445
$ws.Range("L129").Value = 3277.3335
$ws.Range("N129").Value = -13277.3335

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H19").Value = 20933.334
$ws.Range("I19").Value = 2800
$ws.Range("K19").Value = 2800
$ws.Range("M19").Value = -2571

$ws.Range("H33").Value = 5000
$ws.Range("I33").Value = 5000
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 5000
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -4671
$ws.Range("N33").ClearContents()

$ws.Range("H35").Value = 20745.666
$ws.Range("I35").Value = 26118.5
$ws.Range("J35").Value = 10000
$ws.Range("K35").Value = 26118.5
$ws.Range("L35").Value = 10000
$ws.Range("M35").Value = -25712.5
$ws.Range("N35").Value = -10812

$ws.Range("H36").Value = 15731.5
$ws.Range("I36").Value = 14308.667
$ws.Range("K36").Value = 14308.667
$ws.Range("M36").Value = -13962.667

$ws.Range("H76").Value = 29000
$ws.Range("J76").Value = 29000
$ws.Range("L76").Value = 29000
$ws.Range("N76").Value = -29676

$ws.Range("H79").Value = 29000
$ws.Range("J79").Value = 29000
$ws.Range("L79").Value = 29000
$ws.Range("N79").Value = -31340

$ws.Range("H92").Value = 77444.44500000001
$ws.Range("J92").Value = 77444.44500000001
$ws.Range("L92").Value = 77444.44500000001
$ws.Range("N92").Value = -82436.44500000001

$ws.Range("H96").Value = 108705.664
$ws.Range("J96").Value = 108705.664
$ws.Range("L96").Value = 108705.664
$ws.Range("N96").Value = -114197.664

$ws.Range("H102").Value = 1155.4445
$ws.Range("I102").Value = 1155.4445
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1155.4445
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 466.5554999999999
$ws.Range("N102").ClearContents()

$ws.Range("H122").Value = 2091.1875
$ws.Range("I122").Value = 1961.3572
$ws.Range("K122").Value = 5884.071599999999
$ws.Range("M122").Value = -3434.071599999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 14999
$ws.Range("J5").Value = 20000
$ws.Range("L5").Value = 20000
$ws.Range("N5").Value = -20226

$ws.Range("H99").Value = 1956.3125
$ws.Range("I99").Value = 1753.6364
$ws.Range("J99").Value = 2402.2
$ws.Range("K99").Value = 1753.6364
$ws.Range("L99").Value = 2402.2
$ws.Range("M99").Value = -255.6364000000001
$ws.Range("N99").Value = -5398.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H33").Value = 10000
$ws.Range("J33").Value = 10000
$ws.Range("L33").Value = 10000
$ws.Range("N33").Value = -10758

$ws.Range("H132").Value = 1569.7567
$ws.Range("I132").Value = 1320.4286
$ws.Range("J132").Value = 2345.4443
$ws.Range("K132").Value = 3961.2858
$ws.Range("L132").Value = 7036.3329
$ws.Range("M132").Value = -1431.2858
$ws.Range("N132").Value = -12096.3329

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 625618.2
$ws.Range("I113").Value = 2500570.8
$ws.Range("J113").Value = 634
$ws.Range("K113").Value = 7501712.399999999
$ws.Range("L113").Value = 1902
$ws.Range("M113").Value = -7499542.399999999
$ws.Range("N113").Value = -6242

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H88").Value = 49597.5
$ws.Range("J88").Value = 49597.5
$ws.Range("L88").Value = 49597.5
$ws.Range("N88").Value = -50499.5

$ws.Range("H91").Value = 49597.5
$ws.Range("J91").Value = 49597.5
$ws.Range("L91").Value = 49597.5
$ws.Range("N91").Value = -52717.5

$ws.Range("H122").Value = 3125.7334
$ws.Range("I122").Value = 3198.8333
$ws.Range("J122").Value = 2833.3333
$ws.Range("K122").Value = 9596.499899999999
$ws.Range("L122").Value = 8499.999899999999
$ws.Range("M122").Value = -7146.499899999999
$ws.Range("N122").Value = -13399.9999

$ws.Range("H126").Value = 1925.625
$ws.Range("I126").Value = 1650.8334
$ws.Range("J126").Value = 2750
$ws.Range("K126").Value = 4952.5002
$ws.Range("L126").Value = 8250
$ws.Range("M126").Value = -2482.5002
$ws.Range("N126").Value = -13190

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H29").Value = 20000
$ws.Range("J29").Value = 20000
$ws.Range("L29").Value = 20000
$ws.Range("N29").Value = -20590

$ws.Range("H93").Value = 1933.3334
$ws.Range("I93").Value = 1900
$ws.Range("J93").Value = 2000
$ws.Range("K93").Value = 1900
$ws.Range("L93").Value = 2000
$ws.Range("M93").Value = -652
$ws.Range("N93").Value = -4496

$ws.Range("H100").Value = 7550
$ws.Range("I100").Value = 12570
$ws.Range("J100").Value = 3366.6667
$ws.Range("K100").Value = 12570
$ws.Range("L100").Value = 3366.6667
$ws.Range("M100").Value = -12029
$ws.Range("N100").Value = -4448.6667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 49800
$ws.Range("I8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("M8").ClearContents()

$ws.Range("H40").Value = 14333.333
$ws.Range("J40").Value = 14333.333
$ws.Range("L40").Value = 14333.333
$ws.Range("N40").Value = -14631.333

$ws.Range("H64").Value = 22613.5
$ws.Range("J64").Value = 22613.5
$ws.Range("L64").Value = 22613.5
$ws.Range("N64").Value = -23109.5

$ws.Range("H67").Value = 22613.5
$ws.Range("J67").Value = 22613.5
$ws.Range("L67").Value = 22613.5
$ws.Range("N67").Value = -24329.5

$ws.Range("H81").Value = 253250
$ws.Range("I81").Value = 253250
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 506500
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -505439
$ws.Range("N81").ClearContents()

$ws.Range("H84").Value = 253250
$ws.Range("I84").Value = 253250
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 2532500
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -2527196
$ws.Range("N84").ClearContents()
